# Insert a new data row before row 63 (shifts old rows 63..176 down to 64..177)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(63).Insert()

# Fill in the new row 63 with the new record's values (columns A..R).
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across the whole dataset,
# so copy them straight from the row that used to be here (now row 64).
$ws.Cells.Item(63, 1).Value = 10
$ws.Cells.Item(63, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value = "La Araucanía"
$ws.Cells.Item(63, 4).Value = 44757
$ws.Cells.Item(63, 5).Value = 9
$ws.Cells.Item(63, 6).Value = 100114007
$ws.Cells.Item(63, 7).Value = "Jengibre"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 80
$ws.Cells.Item(63, 11).Value = 20000
$ws.Cells.Item(63, 12).Value = 20000
$ws.Cells.Item(63, 13).Value = 20000
$ws.Cells.Item(63, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(63, 15).Value = "Perú"
$ws.Cells.Item(63, 16).Value = 1538
$ws.Cells.Item(63, 17).Value = 13
$ws.Cells.Item(63, 18).Value = "Hortaliza"
